# Updates "Collected Minutiae" minutiae Type column (B2:B36) to the
# corrected w/e values, clears the yellow highlight that flagged the
# previously-wrong cells (now only the genuinely-still-wrong cell stays
# highlighted via its own separate, untouched style), and restores row 10
# to the default row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collected Minutiae")

# Corrected minutiae type per row (shared string "w" or "e").
$values = @{
  2  = "e"
  3  = "e"
  4  = "w"
  5  = "w"
  6  = "w"
  7  = "w"
  8  = "e"
  9  = "w"
  10 = "e"
  11 = "e"
  12 = "e"
  13 = "e"
  14 = "e"
  15 = "e"
  16 = "e"
  17 = "w"
  18 = "w"
  19 = "e"
  20 = "e"
  21 = "e"
  22 = "e"
  23 = "e"
  24 = "e"
  25 = "w"
  26 = "w"
  27 = "e"
  28 = "e"
  29 = "w"
  30 = "e"
  31 = "e"
  32 = "e"
  33 = "e"
  34 = "w"
  35 = "e"
  36 = "w"
}

foreach ($row in $values.Keys) {
  $ws.Range("B$row").Value = $values[$row]
}

# Remove the yellow fill that had marked every row as suspect; the
# remaining genuinely-wrong value keeps its own formatting untouched.
$ws.Range("B2:B36").Interior.Pattern = -4142

# Row 10's height was locked to an old wrapped-text size; reset it to the
# sheet's normal row height.
$ws.Rows.Item(10).RowHeight = 15
